$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency price / 1h-volume figures.
# Every value is written with a leading apostrophe so Excel
# stores it as literal text (these price columns use dotted
# thousands separators, e.g. '69.747.46', which Excel would
# otherwise reinterpret as a number/date). The Style is then
# reset to "Normal" so the quote-prefix does not leave a
# different cell style behind than the original had.

# Row 2
$ws.Range("D2").Value = "'69.747.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.42%  "
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'2.503.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.26%  "
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'576.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.13%  "
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'167.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.44%  "
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.06%  "
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").Value = "'2.503.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.12%  "
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.05%  "
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.97%  "
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("E13").Value = "'  +2.24%  "
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = "'2.948.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.63%  "
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("B15").Value = "'WrappedBTC"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'69.658.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.26%  "
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = "'24.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.18%  "
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = "'2.496.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.51%  "
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.20%  "
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("E20").Value = "'  -3.68%  "
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = "'348.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.63%  "
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'3.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.50%  "
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("E23").Value = "'  +0.65%  "
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("E24").Value = "'  -0.10%  "
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = "'70.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.14%  "
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").Value = "'  +0.35%  "
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'8.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.85%  "
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("D28").Value = "'2.587.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.64%  "
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("E29").Value = "'  -0.57%  "
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("D30").Value = "'0.0₃0893"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.35%  "
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("D31").Value = "'7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.46%  "
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").Value = "'460.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.70%  "
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("E33").Value = "'  -2.95%  "
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("E34").Value = "'  -0.46%  "
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("E35").Value = "'  +0.20%  "
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("E36").Value = "'  +0.33%  "
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("D37").Value = "'157.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.94%  "
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("D38").Value = "'19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.57%  "
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("D39").Value = "'18.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.78%  "
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("E40").Value = "'  +0.03%  "
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.17%  "
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = "'4.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.81%  "
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("E43").Value = "'  +0.36%  "
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'38.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.40%  "
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("E45").Value = "'  -3.78%  "
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("E46").Value = "'  -6.06%  "
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").Value = "'141.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.06%  "
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").Value = "'3.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.29%  "
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").Value = "'0.519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.22%  "
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").Value = "'0.0734"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.56%  "
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = "'0.580"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.53%  "
$ws.Range("E51").Style = "Normal"
